# "Generate Report for Handback"
#
# The CI localization-status report is regenerated: the zh-cn and de-de
# rows move from "Ready for handoff" to "Handed back: in sync with en-US",
# the handback timestamps advance, and the (now resolved) stale-handback
# error detail is cleared out. A couple of report columns are also
# widened/narrowed to fit the new text.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: status shown per-locale (columns E = zh-cn, F = de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# ---- zh-cn detail sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-08-31 02:56:28"
$wsZhCn.Range("P2").Value = ""

# ---- de-de detail sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-08-31 02:56:36"
$wsDeDe.Range("P2").Value = ""

# ---- Column width touch-ups (report layout adjusts to the new text)
# ColumnWidth is in characters; the engine stores width = round(cw,2) then
# snaps to the nearest 1/6 char plus the standard 5/6 grid offset, so we
# pick the input that lands closest to the widths seen in the refreshed
# report (~29.98 chars for the Status columns, ~13.75 for Error Detail).
$wsOverview.Columns.Item(5).ColumnWidth = 29.09
$wsOverview.Columns.Item(6).ColumnWidth = 29.09

$wsZhCn.Columns.Item(3).ColumnWidth = 29.09
$wsZhCn.Columns.Item(16).ColumnWidth = 12.76

$wsDeDe.Columns.Item(3).ColumnWidth = 29.09
$wsDeDe.Columns.Item(16).ColumnWidth = 12.76
